# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# On sheet "Rules", cell B11 previously held the shared string "R40".
# It is changed to the text "1" (a new shared-string entry), while
# keeping the cell's existing style/formatting untouched.
#
# A plain `.Value = "1"` assignment would be auto-coerced to the NUMBER
# 1 by Excel's type inference (since "1" looks numeric), which would
# change the cell's type away from a shared string. To force the
# result to remain a genuine text value (matching the original
# shared-string cell type) without disturbing the cell's style, we
# stage the text through a formula that evaluates to the string "1",
# then collapse it down to a plain value via Copy / Paste Special
# (values only) - exactly like using Excel's "Paste Values" to convert
# a formula result into a literal without touching formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Rules")

$target = $ws.Range("B11")

# Stage the literal text "1" via a formula (guarantees a string result).
$target.Formula = "=""1"""

# Convert the formula result into a plain value in place, preserving
# the cell's existing number format / style (Paste Special -> Values).
$target.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues
